# Tugas Profile Picture management
# Restructure Sheet1: remove the "user_id" column (old column A), shift
# username/nama/level_id contents, and add a new "password" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row
$ws.Range("A1").Value = "level_id"
$ws.Range("B1").Value = "username"
$ws.Range("C1").Value = "nama"
$ws.Range("D1").Value = "password"

# Data rows: level_id, username, nama, password
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "cicawow"
$ws.Range("C2").Value = "caca cici"
$ws.Range("D2").Value = 12345

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "owowow"
$ws.Range("C3").Value = "owowowo"
$ws.Range("D3").Value = 12345

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "kasir"
$ws.Range("C4").Value = "Dela"
$ws.Range("D4").Value = 12345

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "staff"
$ws.Range("C5").Value = "Siska"
$ws.Range("D5").Value = 12345

$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "member"
$ws.Range("C6").Value = "Arif"
$ws.Range("D6").Value = 12345

# Selection moves to E4
$ws.Activate()
$ws.Range("E4").Select()

# Window view size/position change
$win = $excel.ActiveWindow
$win.Left = 1116
$win.Top = 1116
$win.Width = 17280
$win.Height = 8880
